# Apply updated "想去人数" (F column) counts to the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row -> new value for the "展览" sheet (rId1 / sheet1.xml)
$updatesExhibition = @{
    5  = 119
    8  = 72
    9  = 440
    12 = 539
    14 = 278
    16 = 337
    18 = 84
    19 = 46
    20 = 41
    22 = 82
    23 = 829
    24 = 1358
    25 = 283
    26 = 295
    27 = 185
    29 = 151
    31 = 10
    32 = 80
    33 = 197
    35 = 246
    40 = 558
    42 = 3310
    43 = 392
    44 = 172
    45 = 857
}

# Row -> new value for the "全部类型" sheet (rId4 / sheet4.xml)
# Identical to the above except row 42, which differs by one.
$updatesAll = $updatesExhibition.Clone()
$updatesAll[42] = 3311

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $wsExhibition.Range("F$row").Value = $updatesExhibition[$row]
}

$wsAll = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAll.Keys) {
    $wsAll.Range("F$row").Value = $updatesAll[$row]
}
